$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    7  = 90
    8  = 82
    9  = 79
    10 = 88
    11 = 70
    12 = 76
    13 = 72
    14 = 71
    15 = 75
    16 = 80
    17 = 74
    18 = 74
    19 = 71
    20 = 72
    21 = 75
    22 = 68
    23 = 78
    24 = 79
    25 = 69
    26 = 77
    27 = 72
    28 = 72
    29 = 80
    30 = 78
    31 = 84
    33 = 76
    34 = 71
    35 = 66
    36 = 75
    37 = 72
    38 = 68
    39 = 66
    40 = 84
    41 = 68
    42 = 75
    43 = 67
    44 = 68
    45 = 73
    46 = 73
    47 = 70
    48 = 69
    49 = 68
    50 = 73
    51 = 71
    52 = 75
    53 = 69
    54 = 70
    55 = 80
    56 = 68
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
